$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1Bug")

# --- Update header strings (row 1) ---
$ws.Range("D1").Value = "executed_susp_stmt_in_passing_variant"
$ws.Range("F1").Value = "not_executed_susp_stmt_vs_in_passing_variant"
$ws.Range("H1").Value = "executed_susp_stmt_in_a_failed_execution"
$ws.Range("J1").Value = "not_executed_susp_stmt_in_a_failed_execution"

# --- Row 3: /Users/thu-trangnguyen/Documents/Research/SPL/BankAccountTP/1Bug/4wise/ ---
$ws.Range("A3").Value = "/Users/thu-trangnguyen/Documents/Research/SPL/BankAccountTP/1Bug/4wise/"
$ws.Range("B3").Value = 0.68
$ws.Range("C3").Value = 0.87
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0.24
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0.24
$ws.Range("H3").Value = 0.94
$ws.Range("I3").Value = 0.07
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.22
$ws.Range("L3").Value = 0.05
$ws.Range("M3").Value = 0.15
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0.32
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.08
$ws.Range("S3").Value = 0.63
$ws.Range("T3").Value = 0.01
$ws.Range("U3").Value = 0.47
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0.11
$ws.Range("Y3").Value = 0.9
$ws.Range("Z3").Value = 0.19
$ws.Range("AA3").Value = 0.83
$ws.Range("AB3").Value = 0.27
$ws.Range("AC3").Value = 0.76
$ws.Range("AD3").Value = 0.31
$ws.Range("AE3").Value = 0.74
$ws.Range("AF3").Value = 0.3
$ws.Range("AG3").Value = 0.75

# --- Row 4: /Users/thu-trangnguyen/Documents/Research/SPL/Elevator/1Bug/4wise/ ---
$ws.Range("A4").Value = "/Users/thu-trangnguyen/Documents/Research/SPL/Elevator/1Bug/4wise/"
$ws.Range("B4").Value = 0.42
$ws.Range("C4").Value = 0.76
$ws.Range("D4").Value = 0.05
$ws.Range("E4").Value = 0.4
$ws.Range("F4").Value = 0.05
$ws.Range("G4").Value = 0.4
$ws.Range("H4").Value = 0.66
$ws.Range("I4").Value = 0.24
$ws.Range("J4").Value = 0.06
$ws.Range("K4").Value = 0.34
$ws.Range("L4").Value = 0.01
$ws.Range("M4").Value = 0.34
$ws.Range("N4").Value = 0.3
$ws.Range("O4").Value = 0.47
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.32
$ws.Range("S4").Value = 0.62
$ws.Range("T4").Value = 0.11
$ws.Range("U4").Value = 0.42
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 0.22
$ws.Range("Y4").Value = 0.79
$ws.Range("Z4").Value = 0.27
$ws.Range("AA4").Value = 0.79
$ws.Range("AB4").Value = 0.31
$ws.Range("AC4").Value = 0.7
$ws.Range("AD4").Value = 0.26
$ws.Range("AE4").Value = 0.71
$ws.Range("AF4").Value = 0.28
$ws.Range("AG4").Value = 0.71

# --- Row 5: /Users/thu-trangnguyen/Documents/Research/SPL/Email/1Bug/4wise/ ---
$ws.Range("A5").Value = "/Users/thu-trangnguyen/Documents/Research/SPL/Email/1Bug/4wise/"
$ws.Range("B5").Value = 0.5
$ws.Range("C5").Value = 0.55
$ws.Range("D5").Value = 0.01
$ws.Range("E5").Value = 0.38
$ws.Range("F5").Value = 0.01
$ws.Range("G5").Value = 0.38
$ws.Range("H5").Value = 0.84
$ws.Range("I5").Value = 0.29
$ws.Range("J5").Value = 0.01
$ws.Range("K5").Value = 0.34
$ws.Range("L5").Value = 0.04
$ws.Range("M5").Value = 0.35
$ws.Range("N5").Value = 0.05
$ws.Range("O5").Value = 0.43
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0.44
$ws.Range("S5").Value = 0.45
$ws.Range("T5").Value = 0.23
$ws.Range("U5").Value = 0.22
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0.18
$ws.Range("Y5").Value = 0.88
$ws.Range("Z5").Value = 0.34
$ws.Range("AA5").Value = 0.75
$ws.Range("AB5").Value = 0.47
$ws.Range("AC5").Value = 0.81
$ws.Range("AD5").Value = 0.76
$ws.Range("AE5").Value = 0.66
$ws.Range("AF5").Value = 0.57
$ws.Range("AG5").Value = 0.77

# --- Row 6: /Users/thu-trangnguyen/Documents/Research/SPL/ExamDB/1Bug/4wise/ ---
$ws.Range("A6").Value = "/Users/thu-trangnguyen/Documents/Research/SPL/ExamDB/1Bug/4wise/"
$ws.Range("B6").Value = 0.49
$ws.Range("C6").Value = 0.54
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.2
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0.2
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 0.07
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.22
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0.09
$ws.Range("N6").Value = 0.02
$ws.Range("O6").Value = 0.34
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0.44
$ws.Range("S6").Value = 0.81
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0.51
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0.04
$ws.Range("Y6").Value = 0.96
$ws.Range("Z6").Value = 0.08
$ws.Range("AA6").Value = 0.94
$ws.Range("AB6").Value = 0.35
$ws.Range("AC6").Value = 0.75
$ws.Range("AD6").Value = 0.27
$ws.Range("AE6").Value = 0.84
$ws.Range("AF6").Value = 0.33
$ws.Range("AG6").Value = 0.8

# --- Row 7: /Users/thu-trangnguyen/Documents/Research/SPL/GPL/1Bug/1wise/ ---
$ws.Range("A7").Value = "/Users/thu-trangnguyen/Documents/Research/SPL/GPL/1Bug/1wise/"
$ws.Range("B7").Value = 0.45
$ws.Range("C7").Value = 0.69
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.09
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0.09
$ws.Range("H7").Value = 0.86
$ws.Range("I7").Value = 0.08
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.09
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0.29
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0.03
$ws.Range("S7").Value = 0.26
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 0.11
$ws.Range("Y7").Value = 0.92
$ws.Range("Z7").Value = 0.22
$ws.Range("AA7").Value = 0.89
$ws.Range("AB7").Value = 0.25
$ws.Range("AC7").Value = 0.98
$ws.Range("AD7").Value = 0.23
$ws.Range("AE7").Value = 0.95
$ws.Range("AF7").Value = 0.24
$ws.Range("AG7").Value = 0.96

# --- Row 8: /Users/thu-trangnguyen/Documents/Research/SPL/ZipMe/1Bug/2wise/ ---
$ws.Range("A8").Value = "/Users/thu-trangnguyen/Documents/Research/SPL/ZipMe/1Bug/2wise/"
$ws.Range("B8").Value = 0.32
$ws.Range("C8").Value = 0.71
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.69
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0.69
$ws.Range("H8").Value = 0.85
$ws.Range("I8").Value = 0.24
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0.67
$ws.Range("L8").Value = 0.05
$ws.Range("M8").Value = 0.19
$ws.Range("N8").Value = 0.09
$ws.Range("O8").Value = 0.69
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0.34
$ws.Range("S8").Value = 0.63
$ws.Range("T8").Value = 0.04
$ws.Range("U8").Value = 0.46
$ws.Range("V8").Value = 0
$ws.Range("W8").Value = 0
$ws.Range("X8").Value = 0.1
$ws.Range("Y8").Value = 0.92
$ws.Range("Z8").Value = 0.29
$ws.Range("AA8").Value = 0.92
$ws.Range("AB8").Value = 0.52
$ws.Range("AC8").Value = 0.66
$ws.Range("AD8").Value = 0.2
$ws.Range("AE8").Value = 0.68
$ws.Range("AF8").Value = 0.25
$ws.Range("AG8").Value = 0.74

